# Append updated data rows (10-23 August 2021) to the Savignano S.P. report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style reference row (existing row 343) so the new rows match formatting
# of the previous data rows (style "2" on column A = date format).
$styleSourceRow = 343

$newData = @(
    @(344, 44418, 2, 14, 150.7970702283498),
    @(345, 44419, 0, 13, 140.0258509263249),
    @(346, 44420, 1, 11, 118.4834123222749),
    @(347, 44421, 1, 10, 107.7121930202499),
    @(348, 44422, 2, 8,  86.16975441619991),
    @(349, 44423, 3, 9,  96.9409737182249),
    @(350, 44424, 1, 10, 107.7121930202499),
    @(351, 44425, 4, 12, 129.2546316242999),
    @(352, 44426, 0, 12, 129.2546316242999),
    @(353, 44427, 4, 15, 161.5682895303748),
    @(354, 44428, 2, 16, 172.3395088323998),
    @(355, 44429, 0, 14, 150.7970702283498),
    @(356, 44430, 1, 12, 129.2546316242999),
    @(357, 44431, 0, 11, 118.4834123222749)
)

foreach ($entry in $newData) {
    $rowNum = $entry[0]
    $dateSerial = $entry[1]
    $newPos = $entry[2]
    $sumC = $entry[3]
    $sumD = $entry[4]

    # Copy formatting (style) from the last existing data row, then write
    # the new values, matching the layout used by the rest of the sheet.
    $ws.Range("A$styleSourceRow").Copy() | Out-Null
    $ws.Range("A$rowNum").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($rowNum, 1).Value = $dateSerial
    $ws.Cells.Item($rowNum, 2).Value = $newPos
    $ws.Cells.Item($rowNum, 3).Value = $sumC
    $ws.Cells.Item($rowNum, 4).Value = $sumD
}

$excel.CutCopyMode = $false
